$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 (z0bug.invoice_ZI_2 / XE... / 21/TO/1234 / in_invoice)
$ws.Range("F14").Value = "####-<2-99"
$ws.Range("G14").Value = "####-<#-01"
$ws.Range("I14").Value = "####-<#-01"

# Row 15 (z0bug.invoice_ZI_3 / XE125432)
$ws.Range("F15").Value = "####-<#-10"
$ws.Range("G15").Value = "####-<#-10"

# Row 16 (z0bug.invoice_ZI_4 / XE125439) - unchanged

# Row 17 (z0bug.invoice_ZI_5 / FATT/0123/21)
$ws.Range("G17").Value = "####-<#-15"

# Row 18 (z0bug.invoice_ZI_6 / FATT/0124/21)
$ws.Range("G18").Value = "####-<#-20"

# Row 19 (z0bug.invoice_ZI_7 / 21/TO/1590)
$ws.Range("G19").Value = "####-<#-20"

# Row 20 (z0bug.invoice_ZI_8 / TI-8778)
$ws.Range("F20").Value = "####-<#-99"

# Update the saved selection to match the target workbook
$ws.Range("D12").Select()
